# Updated the data with automation users:
# replace the old static "To" / "CC" test values (Anupama D. Thumrugoti /
# Shaik Khaleel) with the new automation account names
# (AutoTestAdmin / AutoTestUser) across all the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transmittals_New")

$ws.Range("A2:A4").Value = "AutoTestAdmin"
$ws.Range("B2:B4").Value = "AutoTestUser"
